# Updates cryptos list table (A1:E51) on Sheet1 with the latest scrape.
# Columns: A=Index(unchanged), B=Coin, C=Link, D=Price, E=Volume(1h)
# D-column values that look numeric are written with a leading "'" so Excel
# keeps them as text (matching the source data, which is all inline/text cells)
# instead of silently parsing them into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "25.770.30"
$ws.Range("E2").Value = "  -0.36%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.583.09"
$ws.Range("E3").Value = "  -2.27%  "

# Row 4: TetherUSD
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.68%  "

# Row 5: BNB
$ws.Range("D5").Value = "'209.13"
$ws.Range("E5").Value = "  -1.69%  "

# Row 6: USDC
$ws.Range("D6").Value = "'0.999"
$ws.Range("E6").Value = "  -0.64%  "

# Row 7: XRP
$ws.Range("D7").Value = "'0.483"
$ws.Range("E7").Value = "  -3.47%  "

# Row 8: Cardano
$ws.Range("D8").Value = "'0.248"
$ws.Range("E8").Value = "  -0.61%  "

# Row 9: Dogecoin
$ws.Range("D9").Value = "'0.0618"
$ws.Range("E9").Value = "  +0.30%  "

# Row 10: Solana
$ws.Range("D10").Value = "'18.22"
$ws.Range("E10").Value = "  -1.32%  "

# Row 11: TRON
$ws.Range("E11").Value = "  -0.52%  "

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.801.28"
$ws.Range("E12").Value = "  -2.36%  "

# Row 13: WrappedEther
$ws.Range("D13").Value = "1.576.07"
$ws.Range("E13").Value = "  -2.70%  "

# Row 14: Polkadot
$ws.Range("D14").Value = "'4.03"
$ws.Range("E14").Value = "  -2.75%  "

# Row 15: Polygon
$ws.Range("D15").Value = "'0.514"
$ws.Range("E15").Value = "  -2.07%  "

# Row 16: WrappedBTC
$ws.Range("D16").Value = "25.758.71"
$ws.Range("E16").Value = "  -0.45%  "

# Row 17: Litecoin
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "'60.26"
$ws.Range("E17").Value = "  -2.00%  "

# Row 18: ShibaInu
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.0₃0723"
$ws.Range("E18").Value = "  -1.81%  "

# Row 19: Dai
$ws.Range("D19").Value = "'0.998"
$ws.Range("E19").Value = "  -0.69%  "

# Row 20: BitcoinCash
$ws.Range("D20").Value = "'191.96"
$ws.Range("E20").Value = "  +0.26%  "

# Row 21: Uniswap
$ws.Range("E21").Value = "  -1.05%  "

# Row 22: Avalanche
$ws.Range("D22").Value = "'9.43"
$ws.Range("E22").Value = "  -0.64%  "

# Row 23: Chainlink
$ws.Range("D23").Value = "'5.94"
$ws.Range("E23").Value = "  -1.44%  "

# Row 24: Stellar
$ws.Range("D24").Value = "'0.131"
$ws.Range("E24").Value = "  -2.81%  "

# Row 25: Monero
$ws.Range("D25").Value = "'141.03"
$ws.Range("E25").Value = "  -1.99%  "

# Row 26: BinanceUSD
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.64%  "

# Row 27: Toncoin
$ws.Range("E27").Value = "  -1.14%  "

# Row 28: EthereumClassic
$ws.Range("D28").Value = "'15.18"
$ws.Range("E28").Value = "  -0.07%  "

# Row 29: Cosmos
$ws.Range("D29").Value = "'6.47"
$ws.Range("E29").Value = "  -2.68%  "

# Row 30: PancakeSwap
$ws.Range("E30").Value = "  -5.67%  "

# Row 31: Hedera
$ws.Range("D31").Value = "'0.0473"
$ws.Range("E31").Value = "  -1.03%  "

# Row 32: Filecoin
$ws.Range("E32").Value = "  -0.30%  "

# Row 33: InternetComputer(DFINITY)
$ws.Range("E33").Value = "  -2.24%  "

# Row 34: LidoDAOToken
$ws.Range("E34").Value = "  +0.80%  "

# Row 35: HuobiToken
$ws.Range("D35").Value = "'2.31"
$ws.Range("E35").Value = "  -4.15%  "

# Row 36: Maker
$ws.Range("D36").Value = "1.101.25"
$ws.Range("E36").Value = "  -2.21%  "

# Row 37: PaxDollar
$ws.Range("B37").Value = "PaxDollar"
$ws.Range("C37").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  -0.71%  "

# Row 38: ImmutableX
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'0.505"
$ws.Range("E38").Value = "  -1.08%  "

# Row 39: VeChain
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0151"
$ws.Range("E39").Value = "  -1.31%  "

# Row 40: MXToken
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.32"
$ws.Range("E40").Value = "  -2.21%  "

# Row 41: ARBITRUM
$ws.Range("D41").Value = "'0.786"
$ws.Range("E41").Value = "  -6.31%  "

# Row 42: TrustWalletToken
$ws.Range("D42").Value = "'0.807"
$ws.Range("E42").Value = "  +7.64%  "

# Row 43: Quant
$ws.Range("D43").Value = "'93.53"
$ws.Range("E43").Value = "  -4.85%  "

# Row 44: FraxShare
$ws.Range("E44").Value = "  +2.27%  "

# Row 45: RocketPoolETH
$ws.Range("D45").Value = "1.713.91"
$ws.Range("E45").Value = "  -2.38%  "

# Row 46: BabyDogeCoin
$ws.Range("D46").Value = "0.0₆0111"
$ws.Range("E46").Value = "  -1.26%  "

# Row 47: RenderToken
$ws.Range("D47").Value = "'1.50"
$ws.Range("E47").Value = "  -0.88%  "

# Row 48: Aave
$ws.Range("D48").Value = "'53.22"
$ws.Range("E48").Value = "  -1.46%  "

# Row 49: Cronos
$ws.Range("E49").Value = "  -1.82%  "

# Row 50: Mantle
$ws.Range("E50").Value = "  -1.37%  "

# Row 51: USDD
$ws.Range("D51").Value = "'0.997"
$ws.Range("E51").Value = "  -0.65%  "
